# Test data for Greece Market
# Adds a new "Greece" worksheet (cloned from "Croatia") at the end of the
# workbook, fills in the Greece-specific market name / ticket reference,
# and moves the "active sheet" focus from Croatia to the new Greece sheet.

$wb = $excel.ActiveWorkbook

# --- Clone the Croatia sheet to the end of the workbook -------------------
$croatia = $wb.Worksheets.Item("Croatia")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Copy([System.Reflection.Missing]::Value, $lastSheet)

$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# --- Fill in the Greece-specific data -------------------------------------
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3167/T3166"

# --- Update selections: Greece becomes the active/selected sheet ----------
$croatia.Select()
$croatia.Cells.Select() | Out-Null

$greece.Select()
$greece.Range("C20").Select() | Out-Null
